$d = $word.ActiveDocument

# --- 1. Paragraph 2: "Why cant I import boostrap and jquery directly from the
#        angular server?" -- split into multiple runs with proofErr markers
#        around the misspelled / flagged words, as a real Word spell-check
#        pass would do. Keep the original pPr (ListParagraph + numPr).
$p2 = $d.Paragraphs(2)
$xml2 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00E40AC9" w:rsidRDefault="00E40AC9" w:rsidP="00E40AC9">' +
        '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
        '<w:r><w:t xml:space="preserve">Why </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:proofErr w:type="gramStart"/>' +
        '<w:r><w:t>cant</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:proofErr w:type="gramEnd"/>' +
        '<w:r><w:t xml:space="preserve"> I import </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>boostrap</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> and </w:t></w:r>' +
        '<w:proofErr w:type="spellStart"/>' +
        '<w:r><w:t>jquery</w:t></w:r>' +
        '<w:proofErr w:type="spellEnd"/>' +
        '<w:r><w:t xml:space="preserve"> directly from the angular server?</w:t></w:r>' +
        '</w:p>'
$p2.Range.InsertXML($xml2)

# --- 2. Paragraph 3: "Ng-CLI doesn't server components, it servers the
#        index.html file" -- drop the _GoBack bookmark that sits in this
#        paragraph (it moves down to the new "Admin account" paragraph
#        below). Keep text + pPr identical.
$p3 = $d.Paragraphs(3)
$xml3 = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidR="00D973DC" w:rsidRPr="00E40AC9" w:rsidRDefault="00D973DC" w:rsidP="00E40AC9">' +
        '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr>' +
        '<w:r><w:t>Ng-CLI doesn' + [char]0x2019 + 't server components, it servers the index.html file</w:t></w:r>' +
        '</w:p>'
$p3.Range.InsertXML($xml3)

# --- 3. Append the new plain paragraphs after paragraph 3, before the
#        section break. Each one is added as a bare paragraph (no style /
#        numbering / rsid inherited), matching the target's plain <w:p>
#        elements. Re-seating the freshly-Added paragraph's content via
#        InsertXML (instead of Range.Text / Style) keeps it free of the
#        w:rsid* cruft Word would normally stamp on a typed paragraph.
function Add-PlainParagraph([string]$text) {
    $np = $d.Paragraphs.Add()
    if ($text -ne $null) {
        $inner = '<w:r><w:t xml:space="preserve">' + $text + '</w:t></w:r>'
    } else {
        $inner = ''
    }
    $xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + $inner + '</w:p>'
    $np.Range.InsertXML($xml)
    return $np
}

# two blank paragraphs
$null = Add-PlainParagraph($null)
$null = Add-PlainParagraph($null)

$null = Add-PlainParagraph("Donation page")
$null = Add-PlainParagraph("Anchor tag")
$null = Add-PlainParagraph("Comment for each letter")
$null = Add-PlainParagraph("Picture carousel")
$null = Add-PlainParagraph("Contact Message")

$pAdmin = Add-PlainParagraph("Admin account")
$d.Bookmarks.Add("_GoBack", $pAdmin.Range)

$null = Add-PlainParagraph("Resume")
$null = Add-PlainParagraph($null)

Write-Output "done"
